$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename header columns (shared strings) ---
$ws.Range("C1").Value = "GDP"
$ws.Range("E1").Value = "Budget_Previous_Year"
$ws.Range("F1").Value = "LatinAmerica"
$ws.Range("G1").Value = "Africa"
$ws.Range("H1").Value = "Confessional"
$ws.Range("I1").Value = "Universal"
$ws.Range("AF1").Value = "Donor_Aid_Budget"
$ws.Range("AG1").Value = "Total_Funds"
$ws.Range("AH1").Value = "%_Private_Funds"
$ws.Range("AI1").Value = "%_MAE_Funds"
$ws.Range("AM1").Value = "Delegation"

# --- Update recalculated numeric values ---
$ws.Range("C2").Value = 2771.04675450926
$ws.Range("C3").Value = 2870.311589353206
$ws.Range("C4").Value = 1873.394108966653
$ws.Range("C5").Value = 1460.056109840828
$ws.Range("C6").Value = 2934.187009790061
$ws.Range("C7").Value = 5191.140356354663
$ws.Range("AL7").Value = 1
$ws.Range("C8").Value = 4729.735976516416
$ws.Range("C9").Value = 4547.50930098406
$ws.Range("C10").Value = 1909.084588129339
$ws.Range("C11").Value = 9502.243585046588
$ws.Range("C12").Value = 6128.19547247793
$ws.Range("C13").Value = 14239.03920301361
$ws.Range("C14").Value = 2100.656463590606
$ws.Range("C15").Value = 1250.795760575873
$ws.Range("C16").Value = 3587.883798243964
$ws.Range("C18").Value = 2839.92516805933
$ws.Range("C19").Value = 2983.242707849043
$ws.Range("C20").Value = 2898.942214704482
$ws.Range("C21").Value = 665.6274194933962
$ws.Range("AL21").Value = 1
$ws.Range("C22").Value = 1904.346464968814
$ws.Range("C23").Value = 1503.870423231357
$ws.Range("C24").Value = 5555.389721901988
$ws.Range("AL24").Value = 1
$ws.Range("C25").Value = 4633.590358399045
$ws.Range("C26").Value = 6336.709213679884
$ws.Range("C27").Value = 1955.461557360978
$ws.Range("C28").Value = 5082.354756663512
$ws.Range("C29").Value = 13825.35808833117
$ws.Range("C30").Value = 2217.474008566157
$ws.Range("C31").Value = 1317.890706178356
$ws.Range("C32").Value = 3579.960081455846
$ws.Range("C34").Value = 2948.84548976845
$ws.Range("C35").Value = 3083.80337578809
$ws.Range("C36").Value = 2965.153206179127
$ws.Range("C37").Value = 691.8942672110555
$ws.Range("AL37").Value = 1
$ws.Range("C38").Value = 1939.33862702996
$ws.Range("C39").Value = 1577.487171555845
$ws.Range("C40").Value = 5660.517066940175
$ws.Range("AL40").Value = 1
$ws.Range("C41").Value = 2024.117324382548
$ws.Range("C42").Value = 6711.616186806423
$ws.Range("C43").Value = 4921.848409120176
$ws.Range("C44").Value = 5360.226632400601
$ws.Range("C45").Value = 2264.394087033834
$ws.Range("C46").Value = 1385.890384668919
$ws.Range("C47").Value = 2094.024217383061
$ws.Range("C48").Value = 5642.578115155247
$ws.Range("C49").Value = 5122.180090208862
$ws.Range("C50").Value = 6911.59200404802
$ws.Range("C51").Value = 5745.422744292303
$ws.Range("AL51").Value = 1
$ws.Range("C52").Value = 1657.651524528445
$ws.Range("C53").Value = 2999.422762626143
$ws.Range("C54").Value = 3156.723844635973
$ws.Range("C55").Value = 1982.009737844954
$ws.Range("C56").Value = 2995.45235738661
$ws.Range("C57").Value = 3748.449444923865
$ws.Range("C59").Value = 2379.668184479739
$ws.Range("C60").Value = 1443.492614888721
$ws.Range("C61").Value = 2201.396847776877
$ws.Range("C62").Value = 5919.20956823756
$ws.Range("C63").Value = 5295.682695961288
$ws.Range("C64").Value = 7200.731056811853
$ws.Range("C65").Value = 5955.175904294275
$ws.Range("AL65").Value = 1
$ws.Range("C66").Value = 1716.389195271215
$ws.Range("C67").Value = 3056.152683606517
$ws.Range("C68").Value = 3212.740625904757
$ws.Range("C69").Value = 2000.792448761861
$ws.Range("C70").Value = 3087.12349650562
$ws.Range("C71").Value = 3796.882621798447
$ws.Range("C73").Value = 2497.68592515536
$ws.Range("C74").Value = 1505.810948829135
$ws.Range("C75").Value = 3843.198240901342
$ws.Range("C76").Value = 2286.013198234259
$ws.Range("C77").Value = 7449.08671983612
$ws.Range("C78").Value = 5412.131646018807
$ws.Range("C79").Value = 3252.634165082374
$ws.Range("C80").Value = 2612.856880840196
$ws.Range("C81").Value = 3137.260298393558
$ws.Range("C82").Value = 2025.814194788851
$ws.Range("C83").Value = 3125.07948072635
$ws.Range("C84").Value = 1775.027517189621
$ws.Range("C85").Value = 5996.49696468919
$ws.Range("C87").Value = 6301.696269820412
$ws.Range("AL87").Value = 1
$ws.Range("C88").Value = 1579.189101937001
$ws.Range("C89").Value = 3748.320622951519
$ws.Range("C90").Value = 2361.056581219794
$ws.Range("C91").Value = 7580.275568826287
$ws.Range("C92").Value = 5330.539154475424
$ws.Range("C93").Value = 3314.741082534716
$ws.Range("C94").Value = 2735.187532014817
$ws.Range("C95").Value = 3210.869677115934
$ws.Range("C96").Value = 2067.29003376698
$ws.Range("C97").Value = 3222.05417836739
$ws.Range("C98").Value = 1836.014008604312
$ws.Range("C99").Value = 6114.227214287786
$ws.Range("C101").Value = 6661.86504232374
$ws.Range("AL101").Value = 1
$ws.Range("C102").Value = 1667.171891046301
$ws.Range("C103").Value = 3530.309422482455
$ws.Range("C104").Value = 2425.561644739583
$ws.Range("C105").Value = 7633.969039669125
$ws.Range("C106").Value = 2854.757682901436
$ws.Range("C107").Value = 5176.058803160127
$ws.Range("C108").Value = 2886.897484630703
$ws.Range("C109").Value = 3242.636921959078
$ws.Range("C110").Value = 3212.81539531051
$ws.Range("C111").Value = 1895.214690888655
$ws.Range("C112").Value = 6262.368904654469
$ws.Range("C113").Value = 7026.178156858586
$ws.Range("AL113").Value = 1
